$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Update existing row 2 values (Candidate ID B2 changes, plus the
# randomized credential columns A/C/D/F/G; E2/H2 keep their values)
$ws.Range("A2").Value = "qyPJB451"
$ws.Range("B2").Value = 231006244
$ws.Range("C2").Value = "mykqzwk73"
$ws.Range("D2").Value = "D79rT#!m"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "EeExZcXb"
$ws.Range("G2").Value = "lxwg"
$ws.Range("H2").Value = "Candidate"

# Add the new row 3 with the same bordered look as row 2
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $cell = $ws.Range($col + "3")
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

$ws.Range("A3").Value = "qZFcP215"
$ws.Range("B3").Value = 231006243
$ws.Range("C3").Value = "xktsvbg93"
$ws.Range("D3").Value = "k24p!SJ%"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "lgbcEkry"
$ws.Range("G3").Value = "tiBC"
$ws.Range("H3").Value = "Candidate"

# Match the widened selection the author left active over the table
$ws.Range("A1:H3").Select() | Out-Null
